$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add recommendation-source (column B) entries for rows 30-33
$ws.Range("B30").Value = "《on writing well》"
$ws.Range("B31").Value = "《详解FPGA：人工智能时代的驱动引擎》"
$ws.Range("B32").Value = "《数字VLSI芯片设计：使用Cadence和Synopsys CAD工具》"
$ws.Range("B33").Value = "《模拟CMOS集成电路设计》"

# Add new book entries to column A, rows 43-47 (row 42 left blank as a separator)
$ws.Range("A43").Value = "《千年一叹》"
$ws.Range("A44").Value = "《明朝那些事儿》"
$ws.Range("A45").Value = "《我与地坛》"
$ws.Range("A46").Value = "《三体》"
$ws.Range("A47").Value = "《查理九世》"

# Last entry added to column B (appended after the new book rows)
$ws.Range("B34").Value = "《SystemVerilog验证：测试平台编写指南》"

# Update the view to match where the user ended up working
$ws.Range("B44").Select()
